$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Sheet ALC
$ws_ALC.Range("H17").Value = 1578.409
$ws_ALC.Range("J17").Value = 1578.409
$ws_ALC.Range("L17").Value = 4735.227000000001
$ws_ALC.Range("N17").Value = -5071.227000000001
$ws_ALC.Range("H20").Value = 4000
$ws_ALC.Range("I20").Value = 4000
$ws_ALC.Range("K20").Value = 4000
$ws_ALC.Range("M20").Value = -3770
$ws_ALC.Range("H35").Value = 4000
$ws_ALC.Range("I35").Value = 4000
$ws_ALC.Range("K35").Value = 4000
$ws_ALC.Range("M35").Value = -3621
$ws_ALC.Range("H51").Value = 7942.5713
$ws_ALC.Range("I51").Value = 2649.75
$ws_ALC.Range("K51").Value = 2649.75
$ws_ALC.Range("M51").Value = -2165.75
$ws_ALC.Range("H69").Value = 9802
$ws_ALC.Range("I69").Value = 8005.5
$ws_ALC.Range("K69").Value = 24016.5
$ws_ALC.Range("M69").Value = -23142.5
$ws_ALC.Range("H72").Value = 9802
$ws_ALC.Range("I72").Value = 8005.5
$ws_ALC.Range("K72").Value = 72049.5
$ws_ALC.Range("M72").Value = -67681.5
$ws_ALC.Range("H74").Value = 4596.5
$ws_ALC.Range("I74").Value = 4596.5
$ws_ALC.Range("J74").Value = 0
$ws_ALC.Range("K74").Value = 4596.5
$ws_ALC.Range("L74").Value = 0
$ws_ALC.Range("M74").Value = -3660.5
$ws_ALC.Range("N74").ClearContents()
$ws_ALC.Range("H77").Value = 4596.5
$ws_ALC.Range("I77").Value = 4596.5
$ws_ALC.Range("J77").Value = 0
$ws_ALC.Range("K77").Value = 22982.5
$ws_ALC.Range("L77").Value = 0
$ws_ALC.Range("M77").Value = -18302.5
$ws_ALC.Range("N77").ClearContents()
$ws_ALC.Range("H80").Value = 3604.1667
$ws_ALC.Range("J80").Value = 3895.6667
$ws_ALC.Range("L80").Value = 11687.0001
$ws_ALC.Range("N80").Value = -13683.0001
$ws_ALC.Range("H83").Value = 3604.1667
$ws_ALC.Range("J83").Value = 3895.6667
$ws_ALC.Range("L83").Value = 35061.0003
$ws_ALC.Range("N83").Value = -45045.0003
$ws_ALC.Range("H103").Value = 522.3333
$ws_ALC.Range("I103").Value = 522.3333
$ws_ALC.Range("J103").Value = 0
$ws_ALC.Range("K103").Value = 1566.9999
$ws_ALC.Range("L103").Value = 0
$ws_ALC.Range("M103").Value = -980.9999
$ws_ALC.Range("N103").ClearContents()
$ws_ALC.Range("H138").Value = 2480.5454
$ws_ALC.Range("I138").Value = 1041.1428
$ws_ALC.Range("J138").Value = 4999.5
$ws_ALC.Range("K138").Value = 3123.4284
$ws_ALC.Range("L138").Value = 14998.5
$ws_ALC.Range("M138").Value = 2016.5716
$ws_ALC.Range("N138").Value = -25278.5

# Sheet ARM
$ws_ARM.Range("H74").Value = 1793.7142
$ws_ARM.Range("I74").Value = 1793.8334
$ws_ARM.Range("K74").Value = 1793.8334
$ws_ARM.Range("M74").Value = -919.8334
$ws_ARM.Range("H77").Value = 1793.7142
$ws_ARM.Range("I77").Value = 1793.8334
$ws_ARM.Range("K77").Value = 8969.166999999999
$ws_ARM.Range("M77").Value = -4601.166999999999

# Sheet BSM
$ws_BSM.Range("H94").Value = 445.66666
$ws_BSM.Range("I94").Value = 304
$ws_BSM.Range("J94").Value = 559
$ws_BSM.Range("K94").Value = 304
$ws_BSM.Range("L94").Value = 559
$ws_BSM.Range("M94").Value = 147
$ws_BSM.Range("N94").Value = -1461
$ws_BSM.Range("H105").Value = 4285.6875
$ws_BSM.Range("I105").Value = 2900.5
$ws_BSM.Range("K105").Value = 2900.5
$ws_BSM.Range("M105").Value = -1153.5
$ws_BSM.Range("H107").Value = 4561.6665
$ws_BSM.Range("I107").Value = 4561.6665
$ws_BSM.Range("K107").Value = 4561.6665
$ws_BSM.Range("M107").Value = -2641.6665

# Sheet CRP
$ws_CRP.Range("H6").Value = 5000
$ws_CRP.Range("I6").Value = 5000
$ws_CRP.Range("K6").Value = 5000
$ws_CRP.Range("M6").Value = -4887
$ws_CRP.Range("H19").Value = 28571700
$ws_CRP.Range("I19").Value = 28571700
$ws_CRP.Range("K19").Value = 28571700
$ws_CRP.Range("M19").Value = -28571530
$ws_CRP.Range("H24").Value = 28571700
$ws_CRP.Range("I24").Value = 28571700
$ws_CRP.Range("K24").Value = 28571700
$ws_CRP.Range("M24").Value = -28571530

# Sheet CUL
$ws_CUL.Range("H51").Value = 0
$ws_CUL.Range("I51").Value = 0
$ws_CUL.Range("K51").Value = 0
$ws_CUL.Range("M51").ClearContents()
$ws_CUL.Range("H131").Value = 2755.7144
$ws_CUL.Range("I131").Value = 2680
$ws_CUL.Range("J131").Value = 2812.5
$ws_CUL.Range("K131").Value = 8040
$ws_CUL.Range("L131").Value = 8437.5
$ws_CUL.Range("M131").Value = -3000
$ws_CUL.Range("N131").Value = -18517.5

# Sheet LTW
$ws_LTW.Range("H22").Value = 1544.3
$ws_LTW.Range("I22").Value = 1075.5
$ws_LTW.Range("J22").Value = 2247.5
$ws_LTW.Range("K22").Value = 1075.5
$ws_LTW.Range("L22").Value = 2247.5
$ws_LTW.Range("M22").Value = -780.5
$ws_LTW.Range("N22").Value = -2837.5
$ws_LTW.Range("H27").Value = 1544.3
$ws_LTW.Range("I27").Value = 1075.5
$ws_LTW.Range("J27").Value = 2247.5
$ws_LTW.Range("K27").Value = 1075.5
$ws_LTW.Range("L27").Value = 2247.5
$ws_LTW.Range("M27").Value = -968.5
$ws_LTW.Range("N27").Value = -2461.5
$ws_LTW.Range("H46").Value = 2828.2222
$ws_LTW.Range("I46").Value = 799.6667
$ws_LTW.Range("J46").Value = 3842.5
$ws_LTW.Range("K46").Value = 799.6667
$ws_LTW.Range("L46").Value = 3842.5
$ws_LTW.Range("M46").Value = -611.6667
$ws_LTW.Range("N46").Value = -4218.5
$ws_LTW.Range("H68").Value = 4732.8887
$ws_LTW.Range("I68").Value = 4732.8887
$ws_LTW.Range("K68").Value = 4732.8887
$ws_LTW.Range("M68").Value = -3983.8887
$ws_LTW.Range("H71").Value = 4732.8887
$ws_LTW.Range("I71").Value = 4732.8887
$ws_LTW.Range("K71").Value = 23664.4435
$ws_LTW.Range("M71").Value = -19920.4435

# Sheet WVR
$ws_WVR.Range("H45").Value = 38990.145
$ws_WVR.Range("I45").Value = 23998
$ws_WVR.Range("K45").Value = 23998
$ws_WVR.Range("M45").Value = -23507
$ws_WVR.Range("H93").Value = 0
$ws_WVR.Range("J93").Value = 0
$ws_WVR.Range("L93").Value = 0
$ws_WVR.Range("N93").ClearContents()
$ws_WVR.Range("H95").Value = 48650
$ws_WVR.Range("J95").Value = 48650
$ws_WVR.Range("L95").Value = 48650
$ws_WVR.Range("N95").Value = -54142
$ws_WVR.Range("H96").Value = 1200
$ws_WVR.Range("I96").Value = 1200
$ws_WVR.Range("K96").Value = 1200
$ws_WVR.Range("M96").Value = 173
$ws_WVR.Range("H97").Value = 24999
$ws_WVR.Range("J97").Value = 24999
$ws_WVR.Range("L97").Value = 24999
$ws_WVR.Range("N97").Value = -26981
$ws_WVR.Range("H100").Value = 1583.3334
$ws_WVR.Range("J100").Value = 4000
$ws_WVR.Range("L100").Value = 8000
$ws_WVR.Range("N100").Value = -9082
$ws_WVR.Range("H130").Value = 36331.668
$ws_WVR.Range("J130").Value = 36331.668
$ws_WVR.Range("L130").Value = 36331.668
$ws_WVR.Range("N130").Value = -46371.668
